# Update the dSF (column F) values for a set of rows to re-pulled data.
# This mirrors the "repull data, push all data, mean calculation" commit,
# which corrected a handful of dSF entries after repulling source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -3
    9  = -10
    10 = -1
    12 = -2
    15 = 7
    17 = -1
    20 = -5
    23 = -3
    27 = -3
    28 = -2
    33 = -5
    34 = -5
    41 = 1
    42 = -7
    44 = -3
    45 = -3
    46 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
